# Auto-generated edit script applying Masamune_Profits market-data updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 30000
$ws.Range("I47").Value = 15000
$ws.Range("J47").Value = 37500
$ws.Range("K47").Value = 15000
$ws.Range("L47").Value = 37500
$ws.Range("M47").Value = -14028
$ws.Range("N47").Value = -39444

$ws.Range("H129").Value = 938.50946
$ws.Range("I129").Value = 1151.6923
$ws.Range("J129").Value = 869.225
$ws.Range("K129").Value = 3455.0769
$ws.Range("L129").Value = 2607.675
$ws.Range("M129").Value = 1544.9231
$ws.Range("N129").Value = -12607.675

$ws.Range("H132").Value = 22248.234
$ws.Range("I132").Value = 2960.205
$ws.Range("J132").Value = 116277.375
$ws.Range("K132").Value = 8880.615
$ws.Range("L132").Value = 348832.125
$ws.Range("M132").Value = -6350.615
$ws.Range("N132").Value = -353892.125

$ws.Range("H137").Value = 7370.927
$ws.Range("I137").Value = 11315.5
$ws.Range("J137").Value = 5738.6895
$ws.Range("K137").Value = 33946.5
$ws.Range("L137").Value = 17216.0685
$ws.Range("M137").Value = -31396.5
$ws.Range("N137").Value = -22316.0685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11856.157
$ws.Range("I32").Value = 11516.122
$ws.Range("K32").Value = 11516.122
$ws.Range("M32").Value = -11229.122

$ws.Range("H44").Value = 38000
$ws.Range("J44").Value = 38000
$ws.Range("L44").Value = 38000
$ws.Range("N44").Value = -38976

$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344

$ws.Range("H132").Value = 20836304
$ws.Range("I132").Value = 33336246
$ws.Range("J132").Value = 3065.7778
$ws.Range("K132").Value = 100008738
$ws.Range("L132").Value = 9197.3334
$ws.Range("M132").Value = -100006208
$ws.Range("N132").Value = -14257.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -860
$ws.Range("N8").ClearContents()

$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990

$ws.Range("H134").Value = 3815
$ws.Range("I134").Value = 4139
$ws.Range("J134").Value = 3410
$ws.Range("K134").Value = 12417
$ws.Range("L134").Value = 10230
$ws.Range("M134").Value = -9882
$ws.Range("N134").Value = -15300

$ws.Range("H135").Value = 69663.19
$ws.Range("J135").Value = 69663.19
$ws.Range("L135").Value = 69663.19
$ws.Range("N135").Value = -79803.19

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 40163620
$ws.Range("I4").Value = 125001310
$ws.Range("J4").Value = 239999.06
$ws.Range("K4").Value = 125001310
$ws.Range("L4").Value = 239999.06
$ws.Range("M4").Value = -125001198
$ws.Range("N4").Value = -240223.06

$ws.Range("H31").Value = 3791228
$ws.Range("I31").Value = 1334.0227
$ws.Range("J31").Value = 7581122
$ws.Range("K31").Value = 1334.0227
$ws.Range("L31").Value = 7581122
$ws.Range("M31").Value = -1039.0227
$ws.Range("N31").Value = -7581712

$ws.Range("H34").Value = 3791228
$ws.Range("I34").Value = 1334.0227
$ws.Range("J34").Value = 7581122
$ws.Range("K34").Value = 1334.0227
$ws.Range("L34").Value = 7581122
$ws.Range("M34").Value = -1132.0227
$ws.Range("N34").Value = -7581526

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H99").Value = 2246.9333
$ws.Range("I99").Value = 2642.8
$ws.Range("J99").Value = 2049
$ws.Range("K99").Value = 2642.8
$ws.Range("L99").Value = 2049
$ws.Range("M99").Value = -1144.8
$ws.Range("N99").Value = -5045

$ws.Range("H126").Value = 2246.9333
$ws.Range("I126").Value = 2642.8
$ws.Range("J126").Value = 2049
$ws.Range("K126").Value = 7928.400000000001
$ws.Range("L126").Value = 6147
$ws.Range("M126").Value = -5458.400000000001
$ws.Range("N126").Value = -11087

$ws.Range("H132").Value = 181038.38
$ws.Range("I132").Value = 2949
$ws.Range("J132").Value = 240401.5
$ws.Range("K132").Value = 8847
$ws.Range("L132").Value = 721204.5
$ws.Range("M132").Value = -6317
$ws.Range("N132").Value = -726264.5

$ws.Range("H134").Value = 610251.9399999999
$ws.Range("I134").Value = 1379.9333
$ws.Range("J134").Value = 1751886.9
$ws.Range("K134").Value = 4139.7999
$ws.Range("L134").Value = 5255660.699999999
$ws.Range("M134").Value = -1604.7999
$ws.Range("N134").Value = -5260730.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 62.5
$ws.Range("I6").Value = 62.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 187.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -74.5
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 15386341
$ws.Range("I113").Value = 20001550
$ws.Range("J113").Value = 2308.6667
$ws.Range("K113").Value = 20001550
$ws.Range("L113").Value = 2308.6667
$ws.Range("M113").Value = -19999380
$ws.Range("N113").Value = -6648.6667

$ws.Range("H132").Value = 3225.1904
$ws.Range("I132").Value = 1935.3334
$ws.Range("J132").Value = 4945
$ws.Range("K132").Value = 5806.0002
$ws.Range("L132").Value = 14835
$ws.Range("M132").Value = -3276.0002
$ws.Range("N132").Value = -19895

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2319.682
$ws.Range("J16").Value = 1998.5
$ws.Range("L16").Value = 1998.5
$ws.Range("N16").Value = -2338.5

$ws.Range("H132").Value = 2623.9167
$ws.Range("I132").Value = 2059.9062
$ws.Range("J132").Value = 3751.9375
$ws.Range("K132").Value = 6179.7186
$ws.Range("L132").Value = 11255.8125
$ws.Range("M132").Value = -3649.7186
$ws.Range("N132").Value = -16315.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 10052502
$ws.Range("I3").Value = 13380001
$ws.Range("K3").Value = 13380001
$ws.Range("M3").Value = -13379887

$ws.Range("H14").Value = 5002250
$ws.Range("J14").Value = 5002250
$ws.Range("L14").Value = 5002250
$ws.Range("N14").Value = -5002586

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H132").Value = 2553.6667
$ws.Range("I132").Value = 2310.7778
$ws.Range("J132").Value = 2699.4
$ws.Range("K132").Value = 6932.3334
$ws.Range("L132").Value = 8098.200000000001
$ws.Range("M132").Value = -4402.3334
$ws.Range("N132").Value = -13158.2
